$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row-209 formatting (style indices for col A / E) into the new rows 210-212 ---
$ws.Range("A209:AC209").Copy()
$ws.Range("A210:AC210").PasteSpecial(-4122)
$ws.Range("A211:AC211").PasteSpecial(-4122)
$ws.Range("A212:AC212").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 209 ---
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = 6787892
$ws.Range("C209").Value = "Croatia HNL"
$ws.Range("D209").Value = "Croatia HNL"
$ws.Range("E209").Value = 45339.45833333334
$ws.Range("F209").Value = "HNK Rijeka"
$ws.Range("G209").Value = "NK Lokomotiva Zagreb"
$ws.Range("H209").Value = 4
$ws.Range("I209").Value = 0
$ws.Range("J209").Value = "H"
$ws.Range("K209").Value = 1.363
$ws.Range("L209").Value = 4.5
$ws.Range("M209").Value = 7
$ws.Range("N209").Value = 1.444
$ws.Range("O209").Value = 4.333
$ws.Range("P209").Value = 6
$ws.Range("Q209").Value = -1.25
$ws.Range("R209").Value = 2.025
$ws.Range("S209").Value = 1.825
$ws.Range("T209").Value = 2.75
$ws.Range("U209").Value = 2
$ws.Range("V209").Value = 1.85
$ws.Range("W209").Value = 0.444
$ws.Range("X209").Value = -1
$ws.Range("Y209").Value = -1
$ws.Range("Z209").Value = 1.025
$ws.Range("AA209").Value = -1
$ws.Range("AB209").Value = 1
$ws.Range("AC209").Value = -1

# --- Row 210 ---
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 6769301
$ws.Range("C210").Value = "Croatia HNL"
$ws.Range("D210").Value = "Croatia HNL"
$ws.Range("E210").Value = 45339.54861111111
$ws.Range("F210").Value = "NK Rudes"
$ws.Range("G210").Value = "Hajduk Split"
$ws.Range("H210").Value = 0
$ws.Range("I210").Value = 2
$ws.Range("J210").Value = "A"
$ws.Range("K210").Value = 7.5
$ws.Range("L210").Value = 4.75
$ws.Range("M210").Value = 1.333
$ws.Range("N210").Value = 10
$ws.Range("O210").Value = 5.5
$ws.Range("P210").Value = 1.25
$ws.Range("Q210").Value = 1.75
$ws.Range("R210").Value = 1.85
$ws.Range("S210").Value = 2
$ws.Range("T210").Value = 2.75
$ws.Range("U210").Value = 1.925
$ws.Range("V210").Value = 1.925
$ws.Range("W210").Value = -1
$ws.Range("X210").Value = -1
$ws.Range("Y210").Value = 0.25
$ws.Range("Z210").Value = -0.5
$ws.Range("AA210").Value = 0.5
$ws.Range("AB210").Value = -1
$ws.Range("AC210").Value = 0.925

# --- Row 211 ---
$ws.Range("A211").Value = 209
$ws.Range("B211").Value = 6788921
$ws.Range("C211").Value = "Croatia HNL"
$ws.Range("D211").Value = "Croatia HNL"
$ws.Range("E211").Value = 45340.45833333334
$ws.Range("F211").Value = "NK Osijek"
$ws.Range("G211").Value = "HNK Gorica"
$ws.Range("H211").Value = 3
$ws.Range("I211").Value = 0
$ws.Range("J211").Value = "H"
$ws.Range("K211").Value = 1.75
$ws.Range("L211").Value = 3.5
$ws.Range("M211").Value = 4.2
$ws.Range("N211").Value = 1.5
$ws.Range("O211").Value = 3.75
$ws.Range("P211").Value = 6.5
$ws.Range("Q211").Value = -1
$ws.Range("R211").Value = 1.85
$ws.Range("S211").Value = 2
$ws.Range("T211").Value = 2.25
$ws.Range("U211").Value = 1.8
$ws.Range("V211").Value = 2.05
$ws.Range("W211").Value = 0.5
$ws.Range("X211").Value = -1
$ws.Range("Y211").Value = -1
$ws.Range("Z211").Value = 0.8500000000000001
$ws.Range("AA211").Value = -1
$ws.Range("AB211").Value = 0.8
$ws.Range("AC211").Value = -1

# --- Row 212 ---
$ws.Range("A212").Value = 210
$ws.Range("B212").Value = 6788920
$ws.Range("C212").Value = "Croatia HNL"
$ws.Range("D212").Value = "Croatia HNL"
$ws.Range("E212").Value = 45340.54861111111
$ws.Range("F212").Value = "Dinamo Zagreb"
$ws.Range("G212").Value = "NK Varazdin"
$ws.Range("H212").Value = 1
$ws.Range("I212").Value = 0
$ws.Range("J212").Value = "H"
$ws.Range("K212").Value = 1.166
$ws.Range("L212").Value = 6
$ws.Range("M212").Value = 15
$ws.Range("N212").Value = 1.25
$ws.Range("O212").Value = 6
$ws.Range("P212").Value = 8.5
$ws.Range("Q212").Value = -1.5
$ws.Range("R212").Value = 1.85
$ws.Range("S212").Value = 2
$ws.Range("T212").Value = 2.75
$ws.Range("U212").Value = 1.975
$ws.Range("V212").Value = 1.875
$ws.Range("W212").Value = 0.25
$ws.Range("X212").Value = -1
$ws.Range("Y212").Value = -1
$ws.Range("Z212").Value = -1
$ws.Range("AA212").Value = 1
$ws.Range("AB212").Value = -1
$ws.Range("AC212").Value = 0.875
